# Since both the scanner and the EEG computer mimic a keyboard, the
# randomisation table's practice participant row (PID 1001) is duplicated
# as a new "dummy" first row (PID 1999) so every real run starts past it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stimuli")

# Insert a new blank row right below the header (row 1), pushing all the
# existing randomisation rows down by one.
$ws.Rows.Item(2).Insert()

# Populate the new row with a duplicate of the very first data row, but
# tag it with PID 1999.
$ws.Cells.Item(2, 1).Value = 1999
$ws.Cells.Item(2, 2).Value = "f1"
$ws.Cells.Item(2, 3).Value = "m1"
$ws.Cells.Item(2, 4).Value = "f3"
$ws.Cells.Item(2, 5).Value = "m3"
$ws.Cells.Item(2, 6).Value = "f4"

# Column D carries a left-border style throughout the table; copy it from
# the row below onto the freshly created cell so formatting matches.
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Drop the stray formatted-but-empty cells that Excel leaves around the
# header/inserted row (H1 originally, H2 from the insert).
$ws.Cells.Item(1, 8).Clear()
$ws.Cells.Item(2, 8).Clear()

# Grow the table ("Table13") so it covers the newly added row.
$lo = $ws.ListObjects.Item("Table13")
$lo.Resize($ws.Range("A1:F43"))

# Make "stimuli" the active sheet/selection (previously "tasks" was
# active), with the cursor left on N13 as in the authored workbook.
$ws.Activate()
$ws.Range("N13").Select()
